$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1867
$ws1.Range("F7").Value = 756
$ws1.Range("F8").Value = 13267
$ws1.Range("F9").Value = 13157
$ws1.Range("F10").Value = 1009
$ws1.Range("F15").Value = 661
$ws1.Range("F16").Value = 2085
$ws1.Range("F18").Value = 42
$ws1.Range("F19").Value = 59
$ws1.Range("F22").Value = 281
$ws1.Range("F24").Value = 10

# Sheet "全部类型" (All types) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1867
$ws4.Range("F9").Value = 756
$ws4.Range("F10").Value = 13267
$ws4.Range("F11").Value = 13157
$ws4.Range("F12").Value = 1009
$ws4.Range("F17").Value = 661
$ws4.Range("F20").Value = 2085
$ws4.Range("F22").Value = 42
$ws4.Range("F23").Value = 59
$ws4.Range("F27").Value = 28
$ws4.Range("F29").Value = 281
$ws4.Range("F33").Value = 10

$wb.Save()
